$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 30
$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 67
$ws.Range("B4").Value = 23
$ws.Range("B5").Value = 54
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 62
$ws.Range("B8").Value = 18
$ws.Range("B9").Value = 88
$ws.Range("B10").Value = 86
